# Apply updated values to column F (dSF) on Sheet1, per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -2
    3  = 0
    4  = 2
    5  = 2
    7  = 1
    8  = -2
    10 = 5
    11 = 0
    12 = -3
    14 = 3
    15 = 1
    17 = -2
    19 = 1
    20 = -1
    21 = 4
    22 = -1
    23 = -1
    24 = -2
    25 = 4
    26 = -1
    27 = 2
    28 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
